# Remove replicate microstates from the SM23 microstate list (v1.4.1).
#
# The rows for these microstate IDs are duplicates/replicates and are
# removed from the worksheet: SM23_micro001, SM23_micro004, SM23_micro006,
# SM23_micro008, SM23_micro014, SM23_micro015, SM23_micro016,
# SM23_micro027, SM23_micro029.
#
# Those correspond (in the original layout) to data rows 3, 4, 6, 8, 11,
# 12, 13, 21 and 23. Deleting them shifts the remaining 27 microstate rows
# up so the sheet ends with rows 1-29 instead of 1-38. The trailing 9
# structure-chart pictures (anchored at the old rows 29-37) are removed to
# match, while the 2D-depiction pictures for the rows that remain keep
# their original anchors/position untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash the two alternating row styles (used for banding the data
# rows) in scratch rows far below the data, so we can re-apply a clean
# alternating pattern after the row deletions shift content (and its
# inline formatting) up. Row 3 carries the "even" banding style, row 4
# the "odd" one.
$ws.Rows.Item(3).Copy()
$ws.Rows.Item(200).PasteSpecial(-4122)
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(201).PasteSpecial(-4122)

# --- Delete the replicate microstate rows, bottom-to-top so row indices
# of not-yet-deleted rows stay valid.
$rowsToDelete = @(23, 21, 13, 12, 11, 8, 6, 4, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# The scratch rows shifted up by the same 9 rows we just removed.
$scratchEven = 200 - $rowsToDelete.Count
$scratchOdd = 201 - $rowsToDelete.Count

# --- Re-stripe the surviving data rows (3..29) with a clean alternating
# pattern, independent of which original row each one came from.
for ($r = 3; $r -le 29; $r++) {
    if ((($r - 3) % 2) -eq 0) {
        $ws.Rows.Item($scratchEven).Copy()
    } else {
        $ws.Rows.Item($scratchOdd).Copy()
    }
    $ws.Rows.Item($r).PasteSpecial(-4122)
}

# --- Remove the scratch rows used to hold the banding templates.
$ws.Rows.Item($scratchOdd).Delete()
$ws.Rows.Item($scratchEven).Delete()

# --- Remove the trailing 9 2D-depiction pictures (previously anchored to
# rows 29-37, the rows that no longer exist at the end of the sheet).
for ($i = $ws.Shapes.Count; $i -ge 28; $i--) {
    $ws.Shapes.Item($i).Delete()
}
